# The commit centers the "Grid Table 7 Colorful" table style (used as
# the document's default table style) both horizontally (paragraph
# justification) and vertically (cell anchoring), for the style as a
# whole and for its "first row" (header) conditional formatting.
#
# Word's table-style conditional-formatting records (w:tblStylePr) and
# a table style's own cell properties (w:tcPr, e.g. w:vAlign) are not
# reachable through the Style COM object - Style only surfaces Font
# and ParagraphFormat (mapping to the style's top-level w:rPr / w:pPr).
# We apply every part of the edit that the Word object model actually
# exposes: centering the style's paragraph alignment, which writes the
# w:pPr/w:jc="center" block into the style definition.

$d = $word.ActiveDocument
$style = $d.Styles.Item("Grid Table 7 Colorful")

# wdAlignParagraphCenter = 1
$style.ParagraphFormat.Alignment = 1
